# Apply updated "dSF" (column F) values as part of a data repull/push.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    7  = 0
    12 = -3
    13 = -4
    18 = 1
    19 = -5
    20 = -1
    22 = 0
    25 = -3
    33 = 5
    34 = -3
    36 = 2
    37 = -3
    39 = -8
    46 = -4
    47 = -2
    51 = -1
    53 = -1
    54 = 2
    55 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
